# Fix errors in template input files:
# rename default file-name placeholders on the CONDUCTOR_files sheet so they
# are prefixed with "template_", widen column E to fit the new text, and
# leave the workbook with CONDUCTOR_COUPLING as the active/selected sheet
# (matching the author's final selection state when they saved the file).

$wb = $excel.ActiveWorkbook

$wsFiles = $wb.Worksheets.Item("CONDUCTOR_files")
$wsInput = $wb.Worksheets.Item("CONDUCTOR_input")
$wsCoupling = $wb.Worksheets.Item("CONDUCTOR_COUPLING")

# --- Update the default file names referenced on CONDUCTOR_files ---
$wsFiles.Range("E4").Value = "template_conductor_1_input.xlsx"
$wsFiles.Range("E5").Value = "template_conductor_1_coupling.xlsx"
$wsFiles.Range("E6").Value = "template_conductor_grid.xlsx"
$wsFiles.Range("E7").Value = "template_conductor_1_operation.xlsx"
$wsFiles.Range("E14").Value = "template_conductor_diagnostic.xlsx"

# Column E needs to be a bit wider to fit the longer file names.
$wsFiles.Columns.Item(5).ColumnWidth = 32.2

# --- Restore sheet selections / active sheet as left by the author ---
$wsFiles.Activate()
$wsFiles.Range("D21").Select() | Out-Null

$wsCoupling.Activate()
$wsCoupling.Range("G17").Select() | Out-Null
